$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove all existing rows/content (rows 1-5) so the sheet starts clean.
$ws.Rows("1:5").Delete()

# --- Row 1: new header row (idx, idx2, Name, Date Start, Date End use the
#     plain/default style; the unit columns F1:K1 use a dedicated style that
#     applies font 1 only, without a number format). ---
$ws.Cells.Item(1,1).Value = "idx"
$ws.Cells.Item(1,2).Value = "idx2"
$ws.Cells.Item(1,3).Value = "Name"
$ws.Cells.Item(1,4).Value = "Date Start"
$ws.Cells.Item(1,5).Value = "Date End"
$ws.Cells.Item(1,6).Value = "(m3/s)"
$ws.Cells.Item(1,7).Value = "(MW1)"
$ws.Cells.Item(1,8).Value = "(MW2)"
$ws.Cells.Item(1,9).Value = "(GWh) Winter"
$ws.Cells.Item(1,10).Value = "(GWh) Summer"
$ws.Cells.Item(1,11).Value = "(GWh) Year"

# Build a temporary named style that reuses the existing "font 1" (Arial 9)
# but applies no explicit number format, then apply it to F1:K1 and remove
# the named style again so the workbook only keeps the resulting direct
# cell format (matching the target cellXfs entry).
$hdrStyle = $wb.Styles.Add("TmpHeaderStyle")
$hdrStyle.Font.Name = "Arial"
$hdrStyle.Font.Size = 9
$ws.Range("F1:K1").Style = "TmpHeaderStyle"
$wb.Styles.Item("TmpHeaderStyle").Delete()

# --- Row 2: Chancy-Pougny ---
$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,1).NumberFormat = "0"
$ws.Cells.Item(2,2).Value = 509700
$ws.Cells.Item(2,2).NumberFormat = "0"
$ws.Cells.Item(2,3).Value = "Chancy-Pougny"
$ws.Cells.Item(2,4).Value = 1925
$ws.Cells.Item(2,4).NumberFormat = "0"
$ws.Cells.Item(2,5).Value = 2008
$ws.Cells.Item(2,5).NumberFormat = "0"
$ws.Cells.Item(2,6).Value = 540
$ws.Cells.Item(2,6).NumberFormat = "0.00"
$ws.Cells.Item(2,7).Value = 31.01
$ws.Cells.Item(2,7).NumberFormat = "0.00"
$ws.Cells.Item(2,8).Value = 28.42
$ws.Cells.Item(2,8).NumberFormat = "0.00"
$ws.Cells.Item(2,9).Value = 64.6
$ws.Cells.Item(2,9).NumberFormat = "0.00"
$ws.Cells.Item(2,10).Value = 83.98
$ws.Cells.Item(2,10).NumberFormat = "0.00"
$ws.Cells.Item(2,11).Value = 148.58
$ws.Cells.Item(2,11).NumberFormat = "0.00"

# --- Row 3: Verbois ---
$ws.Cells.Item(3,1).Value = 2
$ws.Cells.Item(3,1).NumberFormat = "0"
$ws.Cells.Item(3,2).Value = 509600
$ws.Cells.Item(3,2).NumberFormat = "0"
$ws.Cells.Item(3,3).Value = "Verbois"
$ws.Cells.Item(3,4).Value = 1943
$ws.Cells.Item(3,4).NumberFormat = "0"
$ws.Cells.Item(3,5).Value = 1999
$ws.Cells.Item(3,5).NumberFormat = "0"
$ws.Cells.Item(3,6).Value = 620
$ws.Cells.Item(3,6).NumberFormat = "0.00"
$ws.Cells.Item(3,7).Value = 102.8
$ws.Cells.Item(3,7).NumberFormat = "0.00"
$ws.Cells.Item(3,8).Value = 98
$ws.Cells.Item(3,8).NumberFormat = "0.00"
$ws.Cells.Item(3,9).Value = 211
$ws.Cells.Item(3,9).NumberFormat = "0.00"
$ws.Cells.Item(3,10).Value = 255
$ws.Cells.Item(3,10).NumberFormat = "0.00"
$ws.Cells.Item(3,11).Value = 466
$ws.Cells.Item(3,11).NumberFormat = "0.00"

# --- Row 4: Seujet ---
$ws.Cells.Item(4,1).Value = 3
$ws.Cells.Item(4,1).NumberFormat = "0"
$ws.Cells.Item(4,2).Value = 509450
$ws.Cells.Item(4,2).NumberFormat = "0"
$ws.Cells.Item(4,3).Value = "Seujet"
$ws.Cells.Item(4,4).Value = 1994
$ws.Cells.Item(4,4).NumberFormat = "0"
$ws.Cells.Item(4,6).Value = 405
$ws.Cells.Item(4,6).NumberFormat = "0.00"
$ws.Cells.Item(4,7).Value = 8.7
$ws.Cells.Item(4,7).NumberFormat = "0.00"
$ws.Cells.Item(4,8).Value = 5.6
$ws.Cells.Item(4,8).NumberFormat = "0.00"
$ws.Cells.Item(4,9).Value = 9.8
$ws.Cells.Item(4,9).NumberFormat = "0.00"
$ws.Cells.Item(4,10).Value = 10.2
$ws.Cells.Item(4,10).NumberFormat = "0.00"
$ws.Cells.Item(4,11).Value = 20
$ws.Cells.Item(4,11).NumberFormat = "0.00"

# Match the selection recorded in the target workbook.
$ws.Range("A2:K2").Select()
